$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 7500
$ws.Range("I64").Value = 7500
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 7500
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -7252
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 7500
$ws.Range("I67").Value = 7500
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 7500
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -6642
$ws.Range("N67").ClearContents()

$ws.Range("H88").Value = 5725
$ws.Range("I88").Value = 4300
$ws.Range("K88").Value = 4300
$ws.Range("M88").Value = -3894

$ws.Range("H91").Value = 5725
$ws.Range("I91").Value = 4300
$ws.Range("K91").Value = 4300
$ws.Range("M91").Value = -2896

$ws.Range("H100").Value = 2084.25
$ws.Range("J100").Value = 2197.5
$ws.Range("L100").Value = 2197.5
$ws.Range("N100").Value = -3279.5

$ws.Range("H113").Value = 3466
$ws.Range("I113").Value = 3450
$ws.Range("K113").Value = 3450
$ws.Range("M113").Value = -196

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5252.273
$ws.Range("I32").Value = 3567.2222
$ws.Range("K32").Value = 3567.2222
$ws.Range("M32").Value = -3280.2222

$ws.Range("H122").Value = 426286.84
$ws.Range("I122").Value = 505548.25
$ws.Range("K122").Value = 1516644.75
$ws.Range("M122").Value = -1514194.75

$ws.Range("H132").Value = 7998.3335
$ws.Range("J132").Value = 1996
$ws.Range("L132").Value = 5988
$ws.Range("N132").Value = -11048

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3240
$ws.Range("I99").Value = 3240
$ws.Range("K99").Value = 3240
$ws.Range("M99").Value = -1742

$ws.Range("H134").Value = 1124.625
$ws.Range("I134").Value = 1142.4286
$ws.Range("K134").Value = 3427.2858
$ws.Range("M134").Value = -892.2857999999997

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 10090.777
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10090.777
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 10090.777
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -10458.777

$ws.Range("H101").Value = 10090.777
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 10090.777
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 10090.777
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -16580.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1188.3334
$ws.Range("I18").Value = 1188.3334
$ws.Range("K18").Value = 3565.0002
$ws.Range("M18").Value = -3396.0002

$ws.Range("H49").Value = 500
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 500
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 1500
$ws.Range("M49").Value = -1344
$ws.Range("N49").Value = -1812

$ws.Range("H51").Value = 900
$ws.Range("I51").Value = 900
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 2700
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2240
$ws.Range("N51").ClearContents()

$ws.Range("H61").Value = 16.9
$ws.Range("I61").Value = 14.888889
$ws.Range("J61").Value = 35
$ws.Range("K61").Value = 44.666667
$ws.Range("L61").Value = 105
$ws.Range("M61").Value = 170.333333
$ws.Range("N61").Value = -535

$ws.Range("H139").Value = 2799.8
$ws.Range("I139").Value = 2875
$ws.Range("K139").Value = 8625
$ws.Range("M139").Value = -3485

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 2250
$ws.Range("J38").Value = 2250
$ws.Range("L38").Value = 2250
$ws.Range("N38").Value = -3176

$ws.Range("H49").Value = 30000.666
$ws.Range("J49").Value = 30001
$ws.Range("L49").Value = 30001
$ws.Range("N49").Value = -30369

$ws.Range("H92").Value = 12296.333
$ws.Range("J92").Value = 12296.333
$ws.Range("L92").Value = 12296.333
$ws.Range("N92").Value = -16040.333

$ws.Range("H102").Value = 1930.9445
$ws.Range("I102").Value = 1744.3846
$ws.Range("K102").Value = 1744.3846
$ws.Range("M102").Value = -122.3846000000001

$ws.Range("H122").Value = 52056.6
$ws.Range("I122").Value = 1664.625
$ws.Range("K122").Value = 4993.875
$ws.Range("M122").Value = -2543.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 4983.2
$ws.Range("J4").Value = 4003.5
$ws.Range("L4").Value = 4003.5
$ws.Range("N4").Value = -4229.5

$ws.Range("H16").Value = 691.41174
$ws.Range("I16").Value = 691.41174
$ws.Range("K16").Value = 691.41174
$ws.Range("M16").Value = -521.41174

$ws.Range("H22").Value = 1102.8889
$ws.Range("I22").Value = 1146.5714
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 1146.5714
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = -851.5714
$ws.Range("N22").Value = -1540

$ws.Range("H27").Value = 1102.8889
$ws.Range("I27").Value = 1146.5714
$ws.Range("J27").Value = 950
$ws.Range("K27").Value = 1146.5714
$ws.Range("L27").Value = 950
$ws.Range("M27").Value = -1039.5714
$ws.Range("N27").Value = -1164

$ws.Range("H28").Value = 4983.2
$ws.Range("J28").Value = 4003.5
$ws.Range("L28").Value = 4003.5
$ws.Range("N28").Value = -4467.5

$ws.Range("H37").Value = 4983.2
$ws.Range("J37").Value = 4003.5
$ws.Range("L37").Value = 4003.5
$ws.Range("N37").Value = -4217.5

$ws.Range("H46").Value = 2419.9
$ws.Range("J46").Value = 2274.875
$ws.Range("L46").Value = 2274.875
$ws.Range("N46").Value = -2650.875

$ws.Range("H100").Value = 5800.3
$ws.Range("I100").Value = 5875.625
$ws.Range("J100").Value = 5499
$ws.Range("K100").Value = 5875.625
$ws.Range("L100").Value = 5499
$ws.Range("M100").Value = -5334.625
$ws.Range("N100").Value = -6581

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2823.25
$ws.Range("I122").Value = 3179.6
$ws.Range("J122").Value = 2229.3333
$ws.Range("K122").Value = 9538.799999999999
$ws.Range("L122").Value = 6687.999899999999
$ws.Range("M122").Value = -7088.799999999999
$ws.Range("N122").Value = -11587.9999

$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450
